$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 8510
$ws.Cells.Item(62, 9).Value = 9655.714
$ws.Cells.Item(62, 11).Value = 9655.714
$ws.Cells.Item(62, 13).Value = -9031.714
$ws.Cells.Item(65, 8).Value = 8510
$ws.Cells.Item(65, 9).Value = 9655.714
$ws.Cells.Item(65, 11).Value = 48278.57
$ws.Cells.Item(65, 13).Value = -45158.57
$ws.Cells.Item(92, 8).Value = 1316.8572
$ws.Cells.Item(92, 9).Value = 192.88889
$ws.Cells.Item(92, 11).Value = 192.88889
$ws.Cells.Item(92, 13).Value = 1055.11111
$ws.Cells.Item(112, 8).Value = 2682.7812
$ws.Cells.Item(112, 9).Value = 2033.5
$ws.Cells.Item(112, 10).Value = 2977.9092
$ws.Cells.Item(112, 11).Value = 6100.5
$ws.Cells.Item(112, 12).Value = 8933.7276
$ws.Cells.Item(112, 13).Value = -4992.5
$ws.Cells.Item(112, 14).Value = -11149.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 5408
$ws.Cells.Item(25, 9).Value = 816
$ws.Cells.Item(25, 11).Value = 816
$ws.Cells.Item(25, 13).Value = -414
$ws.Cells.Item(26, 8).Value = 1000000
$ws.Cells.Item(26, 9).Value = 1000000
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 1000000
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = -999670
$ws.Cells.Item(26, 14).ClearContents()
$ws.Cells.Item(38, 8).Value = 6314.6665
$ws.Cells.Item(38, 9).Value = 6314.6665
$ws.Cells.Item(38, 11).Value = 6314.6665
$ws.Cells.Item(38, 13).Value = -5847.6665
$ws.Cells.Item(45, 8).Value = 1781.0834
$ws.Cells.Item(45, 9).Value = 1671.75
$ws.Cells.Item(45, 11).Value = 1671.75
$ws.Cells.Item(45, 13).Value = -1294.75
$ws.Cells.Item(133, 8).Value = 74708.664
$ws.Cells.Item(133, 10).Value = 74708.664
$ws.Cells.Item(133, 12).Value = 74708.664
$ws.Cells.Item(133, 14).Value = -79768.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 1439.7693
$ws.Cells.Item(64, 9).Value = 615
$ws.Cells.Item(64, 10).Value = 1955.25
$ws.Cells.Item(64, 11).Value = 615
$ws.Cells.Item(64, 12).Value = 1955.25
$ws.Cells.Item(64, 13).Value = -390
$ws.Cells.Item(64, 14).Value = -2405.25
$ws.Cells.Item(67, 8).Value = 1439.7693
$ws.Cells.Item(67, 9).Value = 615
$ws.Cells.Item(67, 10).Value = 1955.25
$ws.Cells.Item(67, 11).Value = 615
$ws.Cells.Item(67, 12).Value = 1955.25
$ws.Cells.Item(67, 13).Value = 165
$ws.Cells.Item(67, 14).Value = -3515.25
$ws.Cells.Item(99, 8).Value = 1558.4667
$ws.Cells.Item(99, 9).Value = 1444.2222
$ws.Cells.Item(99, 10).Value = 1729.8334
$ws.Cells.Item(99, 11).Value = 1444.2222
$ws.Cells.Item(99, 12).Value = 1729.8334
$ws.Cells.Item(99, 13).Value = 53.77780000000007
$ws.Cells.Item(99, 14).Value = -4725.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 1252500
$ws.Cells.Item(4, 9).Value = 1252500
$ws.Cells.Item(4, 11).Value = 1252500
$ws.Cells.Item(4, 13).Value = -1252388
$ws.Cells.Item(116, 8).Value = 78674
$ws.Cells.Item(116, 9).Value = 78674
$ws.Cells.Item(116, 11).Value = 78674
$ws.Cells.Item(116, 13).Value = -74085

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 18250598
$ws.Cells.Item(4, 9).Value = 27596994
$ws.Cells.Item(4, 11).Value = 82790982
$ws.Cells.Item(4, 13).Value = -82790870
$ws.Cells.Item(132, 8).Value = 1388
$ws.Cells.Item(132, 10).Value = 1565
$ws.Cells.Item(132, 12).Value = 14085
$ws.Cells.Item(132, 14).Value = -19145

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 502833
$ws.Cells.Item(7, 10).Value = 3777.3333
$ws.Cells.Item(7, 12).Value = 3777.3333
$ws.Cells.Item(7, 14).Value = -4001.3333
$ws.Cells.Item(8, 8).Value = 502833
$ws.Cells.Item(8, 10).Value = 3777.3333
$ws.Cells.Item(8, 12).Value = 3777.3333
$ws.Cells.Item(8, 14).Value = -4055.3333
$ws.Cells.Item(122, 8).Value = 46477.832
$ws.Cells.Item(122, 9).Value = 56583
$ws.Cells.Item(122, 10).Value = 8078.2
$ws.Cells.Item(122, 11).Value = 169749
$ws.Cells.Item(122, 12).Value = 24234.6
$ws.Cells.Item(122, 13).Value = -167299
$ws.Cells.Item(122, 14).Value = -29134.6
$ws.Cells.Item(135, 8).Value = 74354
$ws.Cells.Item(135, 10).Value = 74354
$ws.Cells.Item(135, 12).Value = 74354
$ws.Cells.Item(135, 14).Value = -84494

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 7999
$ws.Cells.Item(2, 9).Value = 5998
$ws.Cells.Item(2, 10).Value = 10000
$ws.Cells.Item(2, 11).Value = 5998
$ws.Cells.Item(2, 12).Value = 10000
$ws.Cells.Item(2, 13).Value = -5886
$ws.Cells.Item(2, 14).Value = -10224
$ws.Cells.Item(11, 8).Value = 8000
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 8000
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 8000
$ws.Cells.Item(11, 13).ClearContents()
$ws.Cells.Item(11, 14).Value = -8280
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 6404.294
$ws.Cells.Item(61, 9).Value = 6061.6206
$ws.Cells.Item(61, 10).Value = 8391.799999999999
$ws.Cells.Item(61, 11).Value = 6061.6206
$ws.Cells.Item(61, 12).Value = 8391.799999999999
$ws.Cells.Item(61, 13).Value = -5859.6206
$ws.Cells.Item(61, 14).Value = -8795.799999999999
$ws.Cells.Item(74, 8).Value = 48584.25
$ws.Cells.Item(74, 9).Value = 48584.25
$ws.Cells.Item(74, 11).Value = 48584.25
$ws.Cells.Item(74, 13).Value = -47586.25
$ws.Cells.Item(77, 8).Value = 48584.25
$ws.Cells.Item(77, 9).Value = 48584.25
$ws.Cells.Item(77, 11).Value = 145752.75
$ws.Cells.Item(77, 13).Value = -140760.75
$ws.Cells.Item(82, 8).Value = 1914.2727
$ws.Cells.Item(82, 10).Value = 1181.6666
$ws.Cells.Item(82, 12).Value = 1181.6666
$ws.Cells.Item(82, 14).Value = -1903.6666
$ws.Cells.Item(85, 8).Value = 1914.2727
$ws.Cells.Item(85, 10).Value = 1181.6666
$ws.Cells.Item(85, 12).Value = 1181.6666
$ws.Cells.Item(85, 14).Value = -3677.6666
$ws.Cells.Item(113, 8).Value = 6404.294
$ws.Cells.Item(113, 9).Value = 6061.6206
$ws.Cells.Item(113, 10).Value = 8391.799999999999
$ws.Cells.Item(113, 11).Value = 6061.6206
$ws.Cells.Item(113, 12).Value = 8391.799999999999
$ws.Cells.Item(113, 13).Value = -3891.6206
$ws.Cells.Item(113, 14).Value = -12731.8
$ws.Cells.Item(132, 8).Value = 6494814
$ws.Cells.Item(132, 10).Value = 3846.6667
$ws.Cells.Item(132, 12).Value = 11540.0001
$ws.Cells.Item(132, 14).Value = -16600.0001
$ws.Cells.Item(133, 8).Value = 76459.52
$ws.Cells.Item(133, 10).Value = 67163.336
$ws.Cells.Item(133, 12).Value = 67163.336
$ws.Cells.Item(133, 14).Value = -72223.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 9009.182000000001
$ws.Cells.Item(107, 10).Value = 7139.8
$ws.Cells.Item(107, 12).Value = 21419.4
$ws.Cells.Item(107, 14).Value = -25259.4
$ws.Cells.Item(126, 8).Value = 3952
$ws.Cells.Item(126, 9).Value = 3943.4546
$ws.Cells.Item(126, 11).Value = 11830.3638
$ws.Cells.Item(126, 13).Value = -9360.363799999999
$ws.Cells.Item(139, 8).Value = 78000
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()
